# Refresh the crypto price/volume snapshot in the "cryptos" sheet.
# Values are stored as plain text (inlineStr) in the source workbook, so a
# leading apostrophe is used for any D-column literal that looks like a number
# (e.g. "9.40") to stop Excel from auto-converting it and silently dropping
# trailing zeros / changing its type (the apostrophe itself is not stored).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "68.881.36"
$ws.Range("E2").Value = "  +1.78%  "

# Row 3: Ethereum
$ws.Range("D3").Value = "3.874.00"
$ws.Range("E3").Value = "  +1.04%  "

# Row 4: TetherUSD
$ws.Range("E4").Value = "  +0.10%  "

# Row 5: BNB
$ws.Range("D5").Value = "'602.57"
$ws.Range("E5").Value = "  +1.10%  "

# Row 6: Solana
$ws.Range("D6").Value = "'172.10"
$ws.Range("E6").Value = "  +3.76%  "

# Row 7: LidoStakedEther
$ws.Range("D7").Value = "3.872.10"
$ws.Range("E7").Value = "  +1.13%  "

# Row 8: USDC
$ws.Range("E8").Value = "  +0.11%  "

# Row 9: XRP
$ws.Range("E9").Value = "  +1.13%  "

# Row 10: Dogecoin
$ws.Range("E10").Value = "  +3.32%  "

# Row 11: Toncoin
$ws.Range("D11").Value = "'6.52"
$ws.Range("E11").Value = "  +4.13%  "

# Row 12: Cardano
$ws.Range("D12").Value = "'0.463"
$ws.Range("E12").Value = "  +1.61%  "

# Row 13: ShibaInu
$ws.Range("D13").Value = "'0.0000288"
$ws.Range("E13").Value = "  +16.61%  "

# Row 14: Avalanche
$ws.Range("D14").Value = "'37.31"
$ws.Range("E14").Value = "  +1.50%  "

# Row 15: WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.523.56"
$ws.Range("E15").Value = "  +1.16%  "

# Row 16: WrappedEther
$ws.Range("D16").Value = "3.887.71"
$ws.Range("E16").Value = "  +1.30%  "

# Row 17: WrappedBTC
$ws.Range("D17").Value = "68.835.32"
$ws.Range("E17").Value = "  +1.74%  "

# Row 18: Chainlink
$ws.Range("D18").Value = "'18.47"
$ws.Range("E18").Value = "  +1.50%  "

# Row 19: Polkadot
$ws.Range("E19").Value = "  -0.06%  "

# Row 20: TRON
$ws.Range("E20").Value = "  +0.58%  "

# Row 21: Uniswap
$ws.Range("D21").Value = "'11.16"
$ws.Range("E21").Value = "  +4.33%  "

# Row 22: BitcoinCash
$ws.Range("D22").Value = "'474.28"
$ws.Range("E22").Value = "  +1.54%  "

# Row 23: Polygon
$ws.Range("D23").Value = "'0.732"
$ws.Range("E23").Value = "  +0.74%  "

# Row 24: PEPE
$ws.Range("D24").Value = "'0.0000163"
$ws.Range("E24").Value = "  +3.15%  "

# Row 25: Litecoin
$ws.Range("D25").Value = "'83.91"
$ws.Range("E25").Value = "  +0.69%  "

# Row 26: Fetch.AI
$ws.Range("E26").Value = "  +3.22%  "

# Row 27: InternetComputer(DFINITY)
$ws.Range("E27").Value = "  +0.85%  "

# Row 28: RenderToken
$ws.Range("D28").Value = "'10.51"
$ws.Range("E28").Value = "  +5.50%  "

# Row 29: Dai
$ws.Range("E29").Value = "  +0.23%  "

# Row 30: PancakeSwap
$ws.Range("D30").Value = "'2.96"
$ws.Range("E30").Value = "  +1.73%  "

# Row 31: WrappedeETH
$ws.Range("D31").Value = "4.025.44"
$ws.Range("E31").Value = "  +1.18%  "

# Row 32: NEARProtocol
$ws.Range("D32").Value = "'7.81"
$ws.Range("E32").Value = "  +1.62%  "

# Row 33: EthereumClassic
$ws.Range("D33").Value = "'31.48"
$ws.Range("E33").Value = "  +2.02%  "

# Row 34: ImmutableX
$ws.Range("D34").Value = "'2.31"
$ws.Range("E34").Value = "  +1.11%  "

# Row 35: Aptos
$ws.Range("D35").Value = "'9.40"
$ws.Range("E35").Value = "  +1.11%  "

# Row 36: RenzoRestakedETH
$ws.Range("D36").Value = "3.838.34"
$ws.Range("E36").Value = "  +0.98%  "

# Row 37: dogwifhat
$ws.Range("D37").Value = "'3.97"
$ws.Range("E37").Value = "  +22.21%  "

# Row 38: Hedera
$ws.Range("E38").Value = "  +1.26%  "

# Row 39: Mantle
$ws.Range("E39").Value = "  +1.76%  "

# Row 40: Kaspa
$ws.Range("E40").Value = "  +0.64%  "

# Row 41: Filecoin
$ws.Range("E41").Value = "  +1.95%  "

# Row 42: FirstDigitalUSD
$ws.Range("E42").Value = "  +0.23%  "

# Row 43: TheGraph
$ws.Range("D43").Value = "'0.320"
$ws.Range("E43").Value = "  +2.92%  "

# Row 44: was Stacks, now FLOKI
$ws.Range("B44").Value = "FLOKI"
$ws.Range("C44").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D44").Value = "'0.000303"
$ws.Range("E44").Value = "  +14.21%  "

# Row 45: was FLOKI, now Stacks
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "'2.00"
$ws.Range("E45").Value = "  +1.36%  "

# Row 46: Bittensor
$ws.Range("D46").Value = "'424.23"
$ws.Range("E46").Value = "  +0.48%  "

# Row 48: Cosmos
$ws.Range("D48").Value = "'8.73"
$ws.Range("E48").Value = "  +2.48%  "

# Row 49: OKB
$ws.Range("D49").Value = "'46.44"
$ws.Range("E49").Value = "  -1.53%  "

# Row 50: Monero
$ws.Range("D50").Value = "'142.62"
$ws.Range("E50").Value = "  +0.08%  "

# Row 51: VeChain
$ws.Range("D51").Value = "'0.0359"
$ws.Range("E51").Value = "  +1.50%  "
